$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44518
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("P2").Value = 667

# Row 3
$ws.Range("D3").Value = 44525
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = 8000
$ws.Range("P3").Value = 533

# Row 4
$ws.Range("D4").Value = 44508
$ws.Range("J4").Value = 40
